# Update countries & provincias Spain
# - Swap "Georgia" / "Uruguay" labels at rows 136/137 (order changed in source data)
# - Refresh the "last updated" timestamp string
# - Update the covid-19 stat figures for several countries

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header: update the "last updated" timestamp text (cell A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Junio de 2020 a las 01:36"

# --- Swap the country names for rows 136 and 137 ---
$ws.Range("A136").Value = "Uruguay"
$ws.Range("A137").Value = "Georgia"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 2635501
$ws.Range("C4").Value = 38964
$ws.Range("D4").Value = 1090754
$ws.Range("E4").Value = 1416315
$ws.Range("G4").Value = 280
$ws.Range("H4").Value = 128432

# --- Row 52: Nigeria ---
$ws.Range("B52").Value = 24567
$ws.Range("C52").Value = 490
$ws.Range("D52").Value = 9007
$ws.Range("E52").Value = 14995
$ws.Range("G52").Value = 7
$ws.Range("H52").Value = 565

# --- Row 55: Japon ---
$ws.Range("B55").Value = 18390
$ws.Range("C55").Value = 93
$ws.Range("D55").Value = 16505
$ws.Range("E55").Value = 914

# --- Row 69: Chequia ---
$ws.Range("B69").Value = 11603
$ws.Range("C69").Value = 305
$ws.Range("D69").Value = 7705
$ws.Range("E69").Value = 3550
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 348

# --- Row 72: Noruega ---
$ws.Range("B72").Value = 8855
$ws.Range("C72").Value = 9
$ws.Range("E72").Value = 468

# --- Row 134: Republica de Chipre ---
$ws.Range("D134").Value = 833
$ws.Range("E134").Value = 142

# --- Row 136: now Uruguay (updated figures) ---
$ws.Range("B136").Value = 929
$ws.Range("C136").Value = 5
$ws.Range("D136").Value = 818
$ws.Range("E136").Value = 84
$ws.Range("H136").Value = 27

# --- Row 137: now Georgia (figures formerly shown for Georgia at row 136) ---
$ws.Range("C137").Value = 3
$ws.Range("D137").Value = 785
$ws.Range("E137").Value = 124
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 15

# --- Row 151: Togo ---
$ws.Range("B151").Value = 642
$ws.Range("C151").Value = 27
$ws.Range("D151").Value = 401
$ws.Range("E151").Value = 227

# --- Row 156: Montenegro ---
$ws.Range("E156").Value = 155
$ws.Range("G156").Value = 2
$ws.Range("H156").Value = 11
